# Generate Report for Handoff
# Adds two new file entries (70efb0d1-... and be32e1f2-...) to all three
# worksheets (Overview, zh-cn, de-de) of the localization status report.

$wb = $excel.ActiveWorkbook

# Colors / fonts used by the workbook's "HyperLink" look (blue, underlined).
$hyperlinkColor = 15570276   # BGR encoding of RGB FF6495ED (cornflower blue)
$dateFormat = "yyyy-mm-dd HH:mm:ss"

function Style-AsHyperlink($range) {
    $range.Font.Underline = 2
    $range.Font.Color = $hyperlinkColor
}

function Style-AsDate($range) {
    $range.NumberFormat = $dateFormat
}

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Row 6 - 70efb0d1-c47e-465d-9f32-2fae2ed000da
$wsOverview.Range("A6").Value = "70efb0d1-c47e-465d-9f32-2fae2ed000da.md"
Style-AsHyperlink $wsOverview.Range("A6")
$wsOverview.Hyperlinks.Add($wsOverview.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/c20c2daff9a9ea0e6dab474caad9d9c4f0c8635d/e2e/70efb0d1-c47e-465d-9f32-2fae2ed000da.md", "", "", "70efb0d1-c47e-465d-9f32-2fae2ed000da.md") | Out-Null
$wsOverview.Range("B6").Value = "Ready for handoff"
$wsOverview.Range("C6").Value = "Ready for handoff"
$wsOverview.Range("D6").Value = "2016-03-23 14:41:40"
Style-AsDate $wsOverview.Range("D6")

# Row 7 - be32e1f2-1436-464f-a962-9ec87892f8a8
$wsOverview.Range("A7").Value = "be32e1f2-1436-464f-a962-9ec87892f8a8.md"
Style-AsHyperlink $wsOverview.Range("A7")
$wsOverview.Hyperlinks.Add($wsOverview.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/c20c2daff9a9ea0e6dab474caad9d9c4f0c8635d/e2e/be32e1f2-1436-464f-a962-9ec87892f8a8.md", "", "", "be32e1f2-1436-464f-a962-9ec87892f8a8.md") | Out-Null
$wsOverview.Range("B7").Value = "Ready for handoff"
$wsOverview.Range("C7").Value = "Ready for handoff"
$wsOverview.Range("D7").Value = "2016-03-23 14:41:40"
Style-AsDate $wsOverview.Range("D7")

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Row 6 - 70efb0d1-c47e-465d-9f32-2fae2ed000da
$wsZh.Range("A6").Value = "70efb0d1-c47e-465d-9f32-2fae2ed000da.md"
Style-AsHyperlink $wsZh.Range("A6")
$wsZh.Hyperlinks.Add($wsZh.Range("A6"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/5baed286126a14720f00a435b30230873b41f40e/e2e/70efb0d1-c47e-465d-9f32-2fae2ed000da.md", "", "", "70efb0d1-c47e-465d-9f32-2fae2ed000da.md") | Out-Null
$wsZh.Range("B6").Value = ".md"
$wsZh.Range("C6").Value = "Ready for handoff"
$wsZh.Range("D6").Value = "70efb0d1-c47e-465d-9f32-2fae2ed000da.c1ec72caba1f6f793e43fa952af4c3a40b300a0a.zh-cn.xlf"
Style-AsHyperlink $wsZh.Range("D6")
$wsZh.Hyperlinks.Add($wsZh.Range("D6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c58e12e00f5be8ac2ffe9b748bcd65c06345267a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/70efb0d1-c47e-465d-9f32-2fae2ed000da.c1ec72caba1f6f793e43fa952af4c3a40b300a0a.zh-cn.xlf", "", "", "70efb0d1-c47e-465d-9f32-2fae2ed000da.c1ec72caba1f6f793e43fa952af4c3a40b300a0a.zh-cn.xlf") | Out-Null
$wsZh.Range("E6").Value = "2016-03-23 14:41:36"
Style-AsDate $wsZh.Range("E6")
$wsZh.Range("H6").Value = "0001-01-01 00:00:00"
Style-AsDate $wsZh.Range("H6")
$wsZh.Range("J6").Value = "Include"

# Row 7 - be32e1f2-1436-464f-a962-9ec87892f8a8
$wsZh.Range("A7").Value = "be32e1f2-1436-464f-a962-9ec87892f8a8.md"
Style-AsHyperlink $wsZh.Range("A7")
$wsZh.Hyperlinks.Add($wsZh.Range("A7"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/5baed286126a14720f00a435b30230873b41f40e/e2e/be32e1f2-1436-464f-a962-9ec87892f8a8.md", "", "", "be32e1f2-1436-464f-a962-9ec87892f8a8.md") | Out-Null
$wsZh.Range("B7").Value = ".md"
$wsZh.Range("C7").Value = "Ready for handoff"
$wsZh.Range("D7").Value = "be32e1f2-1436-464f-a962-9ec87892f8a8.2d613c9a9e48981796978e2d6c678ee68e5c09f2.zh-cn.xlf"
Style-AsHyperlink $wsZh.Range("D7")
$wsZh.Hyperlinks.Add($wsZh.Range("D7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c58e12e00f5be8ac2ffe9b748bcd65c06345267a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/be32e1f2-1436-464f-a962-9ec87892f8a8.2d613c9a9e48981796978e2d6c678ee68e5c09f2.zh-cn.xlf", "", "", "be32e1f2-1436-464f-a962-9ec87892f8a8.2d613c9a9e48981796978e2d6c678ee68e5c09f2.zh-cn.xlf") | Out-Null
$wsZh.Range("E7").Value = "2016-03-23 14:41:36"
Style-AsDate $wsZh.Range("E7")
$wsZh.Range("H7").Value = "0001-01-01 00:00:00"
Style-AsDate $wsZh.Range("H7")
$wsZh.Range("J7").Value = "Include"

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# Row 6 - 70efb0d1-c47e-465d-9f32-2fae2ed000da
$wsDe.Range("A6").Value = "70efb0d1-c47e-465d-9f32-2fae2ed000da.md"
Style-AsHyperlink $wsDe.Range("A6")
$wsDe.Hyperlinks.Add($wsDe.Range("A6"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/3a198d039177cce12bbcc256d8d019cf213dea61/e2e/70efb0d1-c47e-465d-9f32-2fae2ed000da.md", "", "", "70efb0d1-c47e-465d-9f32-2fae2ed000da.md") | Out-Null
$wsDe.Range("B6").Value = ".md"
$wsDe.Range("C6").Value = "Ready for handoff"
$wsDe.Range("D6").Value = "70efb0d1-c47e-465d-9f32-2fae2ed000da.c1ec72caba1f6f793e43fa952af4c3a40b300a0a.de-de.xlf"
Style-AsHyperlink $wsDe.Range("D6")
$wsDe.Hyperlinks.Add($wsDe.Range("D6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dd5f7019c4d7f19a25f3bb218ccfc611e938cdd8/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/70efb0d1-c47e-465d-9f32-2fae2ed000da.c1ec72caba1f6f793e43fa952af4c3a40b300a0a.de-de.xlf", "", "", "70efb0d1-c47e-465d-9f32-2fae2ed000da.c1ec72caba1f6f793e43fa952af4c3a40b300a0a.de-de.xlf") | Out-Null
$wsDe.Range("E6").Value = "2016-03-23 14:41:40"
Style-AsDate $wsDe.Range("E6")
$wsDe.Range("H6").Value = "0001-01-01 00:00:00"
Style-AsDate $wsDe.Range("H6")
$wsDe.Range("J6").Value = "Include"

# Row 7 - be32e1f2-1436-464f-a962-9ec87892f8a8
$wsDe.Range("A7").Value = "be32e1f2-1436-464f-a962-9ec87892f8a8.md"
Style-AsHyperlink $wsDe.Range("A7")
$wsDe.Hyperlinks.Add($wsDe.Range("A7"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/3a198d039177cce12bbcc256d8d019cf213dea61/e2e/be32e1f2-1436-464f-a962-9ec87892f8a8.md", "", "", "be32e1f2-1436-464f-a962-9ec87892f8a8.md") | Out-Null
$wsDe.Range("B7").Value = ".md"
$wsDe.Range("C7").Value = "Ready for handoff"
$wsDe.Range("D7").Value = "be32e1f2-1436-464f-a962-9ec87892f8a8.2d613c9a9e48981796978e2d6c678ee68e5c09f2.de-de.xlf"
Style-AsHyperlink $wsDe.Range("D7")
$wsDe.Hyperlinks.Add($wsDe.Range("D7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dd5f7019c4d7f19a25f3bb218ccfc611e938cdd8/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/be32e1f2-1436-464f-a962-9ec87892f8a8.2d613c9a9e48981796978e2d6c678ee68e5c09f2.de-de.xlf", "", "", "be32e1f2-1436-464f-a962-9ec87892f8a8.2d613c9a9e48981796978e2d6c678ee68e5c09f2.de-de.xlf") | Out-Null
$wsDe.Range("E7").Value = "2016-03-23 14:41:40"
Style-AsDate $wsDe.Range("E7")
$wsDe.Range("H7").Value = "0001-01-01 00:00:00"
Style-AsDate $wsDe.Range("H7")
$wsDe.Range("J7").Value = "Include"

Write-Host "Applied handoff report updates for 70efb0d1-c47e-465d-9f32-2fae2ed000da and be32e1f2-1436-464f-a962-9ec87892f8a8"
